$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (2022, Feb) updates
$ws.Range("C13").Value = -4.82559523809525
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = -13.7

# Row 17 (2023, Feb) update
$ws.Range("C17").Value = -6.64583333333334
